# "aggiornamento nomi segment cib"
# Rename the CIB business-line segment labels on the Business_Line sheet:
#   Corporate                 -> Large Corporate - Corporate
#   Business (= Small Business)-> Small Business - SME Retail
#   SME_Corporate              -> SME Corporate (harmonise with the other spelling)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Business_Line")

# Row 3 - "Corporate" segment
$ws.Range("C3").Value = "Large Corporate - Corporate"
$ws.Range("D3").Value = "Large Corporate - Corporate"
$ws.Range("E3").Value = "(SEGMENT = 'Large Corporate - Corporate')"

# Row 4 - "Business" (Small Business) segment
$ws.Range("C4").Value = "Small Business - SME Retail"
$ws.Range("D4").Value = "Small Business - SME Retail"
$ws.Range("E4").Value = "(SEGMENT='Small Business - SME Retail')"

# Row 5 - harmonise "SME_Corporate" with "SME Corporate"
$ws.Range("D5").Value = "SME Corporate"

# Move the active cell / selection on this sheet
$ws.Range("D7").Select()

# Restore the workbook's active sheet / tab selection to "r Threshold_Operator"
$wsThresholdOperator = $wb.Worksheets.Item("r Threshold_Operator")
$wsThresholdOperator.Activate()
$wsThresholdOperator.Range("B15").Select()
